# Civ6 "着力点和时代" workbook update — commit: "[Test Needed] Commemorations and Climate"
#
# Semantic changes (everything else in the raw OOXML diff is Excel-internal
# shared-string/style re-indexing noise caused by re-saving the file, not an
# actual content change):
#   1. Cell A2 label changes from "百花齐放" to "琴棋书画".
#   2. Cell G11 (spy policy text): the red run's trailing cost discount goes
#      from "-90%" to "-75%" — the rest of the rich text (two runs: plain
#      black intro + red body) must keep its original formatting.
#   3. The sheet view's remembered scroll position (topLeftCell="A4") resets
#      back to showing the sheet from the top-left (A1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Plain-text cell — simple value replace is safe (no rich runs here).
$ws.Range("A2").Value = "琴棋书画"

# 2) Rich-text cell — splice just the trailing run's text in place, then
#    restore that run's red font color so the first (plain) run keeps its
#    default/black formatting and the second run keeps its red formatting.
$spyCell = $ws.Range("G11")
$secondRunText = "间谍进行进攻性任务的等级+4。购买间谍需要的成本-75%。"
$secondRunStart = 15
$secondRunLength = $secondRunText.Length

$secondRun = $spyCell.Characters($secondRunStart, $secondRunLength)
$secondRun.Text = $secondRunText

$secondRunAfter = $spyCell.Characters($secondRunStart, $secondRunLength)
$secondRunAfter.Font.Color = 255

# 3) Reset the remembered scroll position so the sheet opens at A1 again
#    (matches topLeftCell="A4" being dropped from the saved view state).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
